$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 90332
$ws.Range("B2").Value = "Bruno Melo"
$ws.Range("C2").Value = "Atendimento ao Cliente"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45084
$ws.Range("G2").Value = 4651.79

# Row 3
$ws.Range("A3").Value = 22565
$ws.Range("B3").Value = "Emilly Freitas"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45101
$ws.Range("G3").Value = 8171.46

# Row 4
$ws.Range("A4").Value = 63601
$ws.Range("B4").Value = "João Vitor Araújo"
$ws.Range("C4").Value = "Operações"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45106
$ws.Range("G4").Value = 2898

# Row 5
$ws.Range("A5").Value = 93166
$ws.Range("B5").Value = "Manuela Santos"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45085
$ws.Range("G5").Value = 3266.84

# Row 6
$ws.Range("A6").Value = 25752
$ws.Range("B6").Value = "Nathan Rodrigues"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45095
$ws.Range("G6").Value = 4843.91

# Row 7
$ws.Range("A7").Value = 58576
$ws.Range("B7").Value = "Davi Luiz Rezende"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45081
$ws.Range("G7").Value = 5627.86

# Row 8
$ws.Range("A8").Value = 57134
$ws.Range("B8").Value = "Srta. Lara da Mota"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45098
$ws.Range("G8").Value = 11389.16

# Row 9
$ws.Range("A9").Value = 86971
$ws.Range("B9").Value = "Sra. Maria Vitória Souza"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45105
$ws.Range("G9").Value = 11250.84

# Row 10
$ws.Range("A10").Value = 85852
$ws.Range("B10").Value = "João Felipe Aragão"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 11997.98

# Row 11
$ws.Range("A11").Value = 86406
$ws.Range("B11").Value = "Giovanna Nascimento"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 4910.56
